# Add working set of sequences: update image/word/category cue data (rows 2-33)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "dog/dog011.jpg"
$ws.Cells.Item(2, 3).Value = "bauen"
$ws.Cells.Item(2, 4).Value = "dog"
$ws.Cells.Item(3, 2).Value = "dog/dog026.jpg"
$ws.Cells.Item(3, 3).Value = "quellen"
$ws.Cells.Item(3, 4).Value = "dog"
$ws.Cells.Item(4, 2).Value = "flower/flower020.jpg"
$ws.Cells.Item(4, 3).Value = "meinen"
$ws.Cells.Item(4, 4).Value = "flower"
$ws.Cells.Item(5, 2).Value = "dog/dog027.jpg"
$ws.Cells.Item(5, 3).Value = "atmen"
$ws.Cells.Item(5, 4).Value = "dog"
$ws.Cells.Item(6, 2).Value = "flower/flower026.jpg"
$ws.Cells.Item(6, 3).Value = "herrschen"
$ws.Cells.Item(6, 4).Value = "flower"
$ws.Cells.Item(7, 2).Value = "dog/dog028.jpg"
$ws.Cells.Item(7, 3).Value = "heben"
$ws.Cells.Item(7, 4).Value = "dog"
$ws.Cells.Item(8, 2).Value = "flower/flower021.jpg"
$ws.Cells.Item(8, 3).Value = "parken"
$ws.Cells.Item(8, 4).Value = "flower"
$ws.Cells.Item(9, 2).Value = "flower/flower002.jpg"
$ws.Cells.Item(9, 3).Value = "rufen"
$ws.Cells.Item(9, 4).Value = "flower"
$ws.Cells.Item(10, 2).Value = "flower/flower001.jpg"
$ws.Cells.Item(10, 3).Value = "fließen"
$ws.Cells.Item(10, 4).Value = "flower"
$ws.Cells.Item(11, 2).Value = "dog/dog007.jpg"
$ws.Cells.Item(11, 3).Value = "kriegen"
$ws.Cells.Item(11, 4).Value = "dog"
$ws.Cells.Item(12, 2).Value = "dog/dog029.jpg"
$ws.Cells.Item(12, 3).Value = "ändern"
$ws.Cells.Item(12, 4).Value = "dog"
$ws.Cells.Item(13, 2).Value = "flower/flower013.jpg"
$ws.Cells.Item(13, 3).Value = "kranken"
$ws.Cells.Item(13, 4).Value = "flower"
$ws.Cells.Item(14, 2).Value = "flower/flower030.jpg"
$ws.Cells.Item(14, 3).Value = "mögen"
$ws.Cells.Item(14, 4).Value = "flower"
$ws.Cells.Item(15, 2).Value = "dog/dog031.jpg"
$ws.Cells.Item(15, 3).Value = "trotzen"
$ws.Cells.Item(15, 4).Value = "dog"
$ws.Cells.Item(16, 2).Value = "flower/flower004.jpg"
$ws.Cells.Item(16, 3).Value = "wecken"
$ws.Cells.Item(16, 4).Value = "flower"
$ws.Cells.Item(17, 2).Value = "flower/flower016.jpg"
$ws.Cells.Item(17, 3).Value = "nullen"
$ws.Cells.Item(17, 4).Value = "flower"
$ws.Cells.Item(18, 2).Value = "dog/dog014.jpg"
$ws.Cells.Item(18, 3).Value = "reisen"
$ws.Cells.Item(18, 4).Value = "dog"
$ws.Cells.Item(19, 2).Value = "dog/dog022.jpg"
$ws.Cells.Item(19, 3).Value = "deuten"
$ws.Cells.Item(19, 4).Value = "dog"
$ws.Cells.Item(20, 2).Value = "dog/dog024.jpg"
$ws.Cells.Item(20, 3).Value = "stoppen"
$ws.Cells.Item(20, 4).Value = "dog"
$ws.Cells.Item(21, 2).Value = "dog/dog018.jpg"
$ws.Cells.Item(21, 3).Value = "küssen"
$ws.Cells.Item(21, 4).Value = "dog"
$ws.Cells.Item(22, 2).Value = "dog/dog012.jpg"
$ws.Cells.Item(22, 3).Value = "spenden"
$ws.Cells.Item(22, 4).Value = "dog"
$ws.Cells.Item(23, 2).Value = "dog/dog010.jpg"
$ws.Cells.Item(23, 3).Value = "münzen"
$ws.Cells.Item(23, 4).Value = "dog"
$ws.Cells.Item(24, 2).Value = "flower/flower019.jpg"
$ws.Cells.Item(24, 3).Value = "piepen"
$ws.Cells.Item(24, 4).Value = "flower"
$ws.Cells.Item(25, 2).Value = "dog/dog015.jpg"
$ws.Cells.Item(25, 3).Value = "narren"
$ws.Cells.Item(25, 4).Value = "dog"
$ws.Cells.Item(26, 2).Value = "dog/dog009.jpg"
$ws.Cells.Item(26, 3).Value = "binden"
$ws.Cells.Item(26, 4).Value = "dog"
$ws.Cells.Item(27, 2).Value = "dog/dog021.jpg"
$ws.Cells.Item(27, 3).Value = "streifen"
$ws.Cells.Item(27, 4).Value = "dog"
$ws.Cells.Item(28, 2).Value = "flower/flower009.jpg"
$ws.Cells.Item(28, 3).Value = "grenzen"
$ws.Cells.Item(28, 4).Value = "flower"
$ws.Cells.Item(29, 2).Value = "flower/flower022.jpg"
$ws.Cells.Item(29, 3).Value = "kennen"
$ws.Cells.Item(29, 4).Value = "flower"
$ws.Cells.Item(30, 2).Value = "flower/flower005.jpg"
$ws.Cells.Item(30, 3).Value = "legen"
$ws.Cells.Item(30, 4).Value = "flower"
$ws.Cells.Item(31, 2).Value = "flower/flower023.jpg"
$ws.Cells.Item(31, 3).Value = "tollen"
$ws.Cells.Item(31, 4).Value = "flower"
$ws.Cells.Item(32, 2).Value = "flower/flower028.jpg"
$ws.Cells.Item(32, 3).Value = "wehen"
$ws.Cells.Item(32, 4).Value = "flower"
$ws.Cells.Item(33, 2).Value = "flower/flower029.jpg"
$ws.Cells.Item(33, 3).Value = "passen"
$ws.Cells.Item(33, 4).Value = "flower"
